$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 510.2
$ws.Range("I11").Value = 510.2
$ws.Range("K11").Value = 510.2
$ws.Range("M11").Value = -370.2
$ws.Range("H17").Value = 1999.3334
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1999.3334
$ws.Range("K17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("M17").Value = 5998.0002
$ws.Range("N17").Value = -6334.0002
$ws.Range("H33").Value = 266.17648
$ws.Range("J33").Value = 1049
$ws.Range("L33").Value = 1049
$ws.Range("N33").Value = -1507
$ws.Range("H76").Value = 6824.5625
$ws.Range("I76").Value = 6028
$ws.Range("J76").Value = 7444.1113
$ws.Range("K76").Value = 6028
$ws.Range("L76").Value = 7444.1113
$ws.Range("M76").Value = -5713
$ws.Range("N76").Value = -8074.1113
$ws.Range("H79").Value = 6824.5625
$ws.Range("I79").Value = 6028
$ws.Range("J79").Value = 7444.1113
$ws.Range("K79").Value = 6028
$ws.Range("L79").Value = 7444.1113
$ws.Range("M79").Value = -4936
$ws.Range("N79").Value = -9628.1113
$ws.Range("H87").Value = 80000.5
$ws.Range("J87").Value = 80000.5
$ws.Range("L87").Value = 80000.5
$ws.Range("N87").Value = -82496.5
$ws.Range("H90").Value = 80000.5
$ws.Range("J90").Value = 80000.5
$ws.Range("L90").Value = 240001.5
$ws.Range("N90").Value = -252481.5
$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -44992
$ws.Range("H95").Value = 26999
$ws.Range("J95").Value = 26999
$ws.Range("L95").Value = 26999
$ws.Range("N95").Value = -32491
$ws.Range("H137").Value = 2471.7407
$ws.Range("I137").Value = 1295.2
$ws.Range("J137").Value = 3163.8235
$ws.Range("K137").Value = 3885.6
$ws.Range("L137").Value = 9491.470499999999
$ws.Range("M137").Value = -1335.6
$ws.Range("N137").Value = -14591.4705
$ws.Range("H138").Value = 6963.759
$ws.Range("I138").Value = 4414.4
$ws.Range("K138").Value = 13243.2
$ws.Range("M138").Value = -8103.199999999999
$ws.Range("H141").Value = 8161
$ws.Range("I141").Value = 8161
$ws.Range("K141").Value = 24483
$ws.Range("M141").Value = -19303

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1959
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H63").Value = 8722.223
$ws.Range("I63").Value = 7125
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 7125
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -6439
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 8722.223
$ws.Range("I66").Value = 7125
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 35625
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -32193
$ws.Range("N66").Value = -56864
$ws.Range("H74").Value = 2134.25
$ws.Range("I74").Value = 896.7
$ws.Range("J74").Value = 4196.8335
$ws.Range("K74").Value = 896.7
$ws.Range("L74").Value = 4196.8335
$ws.Range("M74").Value = -22.70000000000005
$ws.Range("N74").Value = -5944.8335
$ws.Range("H77").Value = 2134.25
$ws.Range("I77").Value = 896.7
$ws.Range("J77").Value = 4196.8335
$ws.Range("K77").Value = 4483.5
$ws.Range("L77").Value = 20984.1675
$ws.Range("M77").Value = -115.5
$ws.Range("N77").Value = -29720.1675
$ws.Range("H136").Value = 1959
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2604.75
$ws.Range("I86").Value = 1204
$ws.Range("K86").Value = 1204
$ws.Range("M86").Value = -81
$ws.Range("H89").Value = 2604.75
$ws.Range("I89").Value = 1204
$ws.Range("K89").Value = 6020
$ws.Range("M89").Value = -404
$ws.Range("H99").Value = 1412
$ws.Range("I99").Value = 1412
$ws.Range("K99").Value = 1412
$ws.Range("M99").Value = 86
$ws.Range("H105").Value = 4613.5386
$ws.Range("I105").Value = 3872.0625
$ws.Range("J105").Value = 5799.9
$ws.Range("K105").Value = 3872.0625
$ws.Range("L105").Value = 5799.9
$ws.Range("M105").Value = -2125.0625
$ws.Range("N105").Value = -9293.9

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3897.625
$ws.Range("I31").Value = 2185.7778
$ws.Range("K31").Value = 2185.7778
$ws.Range("M31").Value = -1890.7778
$ws.Range("H34").Value = 3897.625
$ws.Range("I34").Value = 2185.7778
$ws.Range("K34").Value = 2185.7778
$ws.Range("M34").Value = -1983.7778
$ws.Range("H58").Value = 5076
$ws.Range("J58").Value = 6345.1665
$ws.Range("L58").Value = 6345.1665
$ws.Range("N58").Value = -6751.1665
$ws.Range("H107").Value = 475.06668
$ws.Range("I107").Value = 379.53845
$ws.Range("K107").Value = 379.53845
$ws.Range("M107").Value = 1540.46155
$ws.Range("H136").Value = 5076
$ws.Range("J136").Value = 6345.1665
$ws.Range("L136").Value = 19035.4995
$ws.Range("N136").Value = -24135.4995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 52.076923
$ws.Range("J2").Value = 95.8
$ws.Range("L2").Value = 574.8
$ws.Range("N2").Value = -800.8
$ws.Range("H32").Value = 3671295.8
$ws.Range("I32").Value = 2998.6667
$ws.Range("J32").Value = 5364356
$ws.Range("K32").Value = 8996.000100000001
$ws.Range("L32").Value = 16093068
$ws.Range("M32").Value = -8713.000100000001
$ws.Range("N32").Value = -16093634
$ws.Range("H81").Value = 1957.5
$ws.Range("J81").Value = 1957.5
$ws.Range("L81").Value = 5872.5
$ws.Range("N81").Value = -8118.5
$ws.Range("H84").Value = 1957.5
$ws.Range("J84").Value = 1957.5
$ws.Range("L84").Value = 17617.5
$ws.Range("N84").Value = -28849.5
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H121").Value = 992.5
$ws.Range("I121").Value = 296.66666
$ws.Range("J121").Value = 1410
$ws.Range("K121").Value = 889.9999799999999
$ws.Range("L121").Value = 4230
$ws.Range("M121").Value = 420.0000200000001
$ws.Range("N121").Value = -6850

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 9999.333000000001
$ws.Range("J33").Value = 9999.333000000001
$ws.Range("L33").Value = 9999.333000000001
$ws.Range("N33").Value = -10503.333
$ws.Range("H102").Value = 1221.3684
$ws.Range("I102").Value = 648.7083
$ws.Range("K102").Value = 648.7083
$ws.Range("M102").Value = 973.2917
$ws.Range("H107").Value = 908.5
$ws.Range("I107").Value = 728
$ws.Range("K107").Value = 728
$ws.Range("M107").Value = 1192
$ws.Range("H122").Value = 613811.9399999999
$ws.Range("I122").Value = 78817.62
$ws.Range("K122").Value = 236452.86
$ws.Range("M122").Value = -234002.86

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1925.7
$ws.Range("I40").Value = 1657.125
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 1657.125
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1521.125
$ws.Range("N40").Value = -3272
$ws.Range("H122").Value = 3550.7778
$ws.Range("I122").Value = 3524.5715
$ws.Range("K122").Value = 10573.7145
$ws.Range("M122").Value = -8123.7145
$ws.Range("H132").Value = 5399.636
$ws.Range("I132").Value = 3032
$ws.Range("J132").Value = 6287.5
$ws.Range("K132").Value = 9096
$ws.Range("L132").Value = 18862.5
$ws.Range("M132").Value = -6566
$ws.Range("N132").Value = -23922.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1338330.6
$ws.Range("J5").Value = 7496
$ws.Range("L5").Value = 7496
$ws.Range("N5").Value = -7720
$ws.Range("H113").Value = 2574.7646
$ws.Range("I113").Value = 2235
$ws.Range("K113").Value = 6705
$ws.Range("M113").Value = -4535
$ws.Range("H126").Value = 73849.28999999999
$ws.Range("I126").Value = 101239.5
$ws.Range("K126").Value = 303718.5
$ws.Range("M126").Value = -301248.5
$ws.Range("H132").Value = 2100.074
$ws.Range("I132").Value = 1700.1904
$ws.Range("K132").Value = 5100.5712
$ws.Range("M132").Value = -2570.5712
$ws.Range("H136").Value = 25748.56
$ws.Range("I136").Value = 1276.08
$ws.Range("J136").Value = 63986.812
$ws.Range("K136").Value = 3828.24
$ws.Range("L136").Value = 191960.436
$ws.Range("M136").Value = -1278.24
$ws.Range("N136").Value = -197060.436
